$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)

# Update cell values
$ws.Range("A2").Value = "ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("A3").Value = "ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("A4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# Rebuild hyperlinks in the new row order (identity-stable URLs preserved)
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md", [Type]::Missing, [Type]::Missing, "ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md", [Type]::Missing, [Type]::Missing, "ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/39def177-56cc-4c9f-a38f-5fa62a4c38ed.md", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)

# Update cell values
$ws.Range("A2").Value = "ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-08 08:41:08"
$ws.Range("E2").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.md"
$ws.Range("F2").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-08 08:41:27"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-08 08:41:08"
$ws.Range("E3").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.md"
$ws.Range("F3").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf"
$ws.Range("G3").Value = "2016-03-08 08:41:27"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-08 08:43:43"
$ws.Range("E4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md"
$ws.Range("F4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.zh-cn.xlf"
$ws.Range("G4").Value = "2016-03-08 08:43:16"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks in the new row order (identity-stable URLs preserved)
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md", [Type]::Missing, [Type]::Missing, "ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/913ae773c1c985759f13b1c397a9cb38444880b1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f30eb279fba10d916c0d89868937c9e469edf280/e2e/13ce8274-55a2-443c-bc8a-e01f706c76e3.md", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1b39fe16cc28a2081722dceefd2baed96c3289ad/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md", [Type]::Missing, [Type]::Missing, "ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/913ae773c1c985759f13b1c397a9cb38444880b1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f30eb279fba10d916c0d89868937c9e469edf280/e2e/13ce8274-55a2-443c-bc8a-e01f706c76e3.md", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1b39fe16cc28a2081722dceefd2baed96c3289ad/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/39def177-56cc-4c9f-a38f-5fa62a4c38ed.md", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4f282a01e29116718595cc27efd546a7cd8c584f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3fa40997c991c1bae6e4d6972da15faaa1b7074e/e2e/39def177-56cc-4c9f-a38f-5fa62a4c38ed.md", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e897d5f26d8a4f41ee45a363d207e4bc2e457421/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)

# Update cell values
$ws.Range("A2").Value = "ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf"
$ws.Range("D2").Value = "2016-03-08 08:41:12"
$ws.Range("E2").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.md"
$ws.Range("F2").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf"
$ws.Range("G2").Value = "2016-03-08 08:41:34"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf"
$ws.Range("D3").Value = "2016-03-08 08:41:12"
$ws.Range("E3").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.md"
$ws.Range("F3").Value = "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf"
$ws.Range("G3").Value = "2016-03-08 08:41:34"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.de-de.xlf"
$ws.Range("D4").Value = "2016-03-08 08:43:47"
$ws.Range("E4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md"
$ws.Range("F4").Value = "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.de-de.xlf"
$ws.Range("G4").Value = "2016-03-08 08:43:24"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks in the new row order (identity-stable URLs preserved)
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md", [Type]::Missing, [Type]::Missing, "ffff9887e4a4-bf3b-4507-bec4-b3abcd65115e.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/857913ee9d25a1f7887fe9b1ab5856180598805d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4117ca40259154c511547b24a912dedca8228882/e2e/13ce8274-55a2-443c-bc8a-e01f706c76e3.md", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c8ea08d09d0f49d0484520c77194d36f72db1764/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md", [Type]::Missing, [Type]::Missing, "ffffff9c9bdc51-7168-4818-a6a6-874aa336f7c7.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/857913ee9d25a1f7887fe9b1ab5856180598805d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4117ca40259154c511547b24a912dedca8228882/e2e/13ce8274-55a2-443c-bc8a-e01f706c76e3.md", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c8ea08d09d0f49d0484520c77194d36f72db1764/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf", [Type]::Missing, [Type]::Missing, "13ce8274-55a2-443c-bc8a-e01f706c76e3.ccba7915c610d7dd347a6c0e7c036de7f5ecc82c.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/e2e/39def177-56cc-4c9f-a38f-5fa62a4c38ed.md", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/153ef42a01632c04a8b7c698252761dabb235eb0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.de-de.xlf", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ddf62fdf8a64b669a5d0ee174768ad8bdbef62a6/e2e/39def177-56cc-4c9f-a38f-5fa62a4c38ed.md", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.md")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f8cf2bf1389aba059caf8bd954436f77962c6a49/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.de-de.xlf", [Type]::Missing, [Type]::Missing, "39def177-56cc-4c9f-a38f-5fa62a4c38ed.cdecb388ca1ec1719d412a7621e90b538a5dafd8.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c920f232be352f73e489cdf480f0b136e4238bc5/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")
